# ----------------------------------------------------------------------------
# Commit: "add graphAnalysis LoadData ,RenameTab,QueryLoadStatus case"
#
# 1. Resize the saved workbook window (best-effort; some COM hosts only
#    persist Width, not Height, of the window chrome).
# 2. Add a new "isRun" column (E) to the loadData test-case sheet, with a
#    flag value per test row (1 for the two "should run" rows, 0 for the
#    rest).
# 3. Parametrize the first test case's JSON params so projectId/graphId are
#    placeholders (${projectId} / ${graphId}) instead of hard-coded values.
# 4. Move the active selection to F4 (matches the author's last cursor spot
#    when they saved).
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. window size -----------------------------------------------------
$win = $excel.ActiveWindow
$win.Width  = 28000
$win.Height = 13120

# --- 2. new "isRun" header + flag column --------------------------------
$ws.Range("E1").Value = "isRun"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

# --- 3. parametrize projectId / graphId in the first test case ----------
$ws.Range("B2").Value = '{"projectId":${projectId},"graphId":${graphId},"fileName":"graphTestData.json"}'

# --- 4. move selection to F4 --------------------------------------------
$ws.Range("F4").Select()
